# Update cryptocurrency price/volume snapshot values (GitHub Actions scrape refresh).
# Price-column values that read as plain decimals are forced through a text
# NumberFormat so Excel stores them as literal strings (matching the source feed's
# formatting, e.g. trailing zeros) instead of silently parsing them as numbers;
# ClearFormats() afterwards drops the temporary format so no style residue is left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.165.38"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "3.413.82"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +6.92%  "
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.24"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000281"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "678.46"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.05%  "
$ws.Range("D14").Value = "3.962.12"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.63"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").Value = "69.310.63"
$ws.Range("D17").Value = "3.419.22"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  -2.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.01"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.76"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.64"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.49"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.73"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.86"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.70"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.62%  "
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "550.40"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.02"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").Value = "3.608.55"
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").Value = "0.0₃0734"
$ws.Range("E40").Value = "  +8.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.26"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.68"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.39"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("E44").Value = "  +2.73%  "
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  +4.58%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.63"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.24%  "
